# Fix "num hrs late" computation for the DTR entry in row 5 (02-10-2015):
# it was missing its "NO OF HOURS LATE" value and was still carrying the
# highlighted ("needs review") row style left over from data entry.
# Bring the row's formatting back to the normal data-row style (same as
# row 6, etc.) and fill in the late-hours figure (0.25 = 15 minutes late).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the normal row style (no fill highlight) from row 6 onto row 5,
# reusing the existing "normal" cell format instead of creating a new one.
$ws.Range("A6:J6").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Set the "NO OF HOURS LATE" (column F) value for this row.
$ws.Range("F5").Value = 0.25

Write-Host "done"
